$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "Last status check on: 30.01.2022 16:45"

$ws.Range("D10").Value = 0.2
$ws.Range("E10").Value = 44591.69016203703
$ws.Range("E10").NumberFormat = "YYYY-MM-DD HH:MM:SS"
